# "Las funciones de dano funcionan correctamente con marco de un nivel"
#
# This script extends several input sheets of the "Marco 3D" workbook with
# extra node/connectivity/support data for a one-level frame (nodes 5-12,
# elements 4-12) and fills in the "prop geom" helper rows (4-12) that feed
# the damage-function formulas on the "opensees" sheet, so every downstream
# formula (N8/N9/P8/Q8/R8/... on "opensees") recalculates to non-zero values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "nudos" (node coordinates) - add nodes 5-8, fix D1, grow selection
# ---------------------------------------------------------------------
$nudos = $wb.Worksheets.Item("nudos")

$nudos.Cells.Item(1,4).Value = 5000

$nudos.Cells.Item(5,1).Value = 5
$nudos.Cells.Item(5,2).Value = 0
$nudos.Cells.Item(5,3).Value = 5000
$nudos.Cells.Item(5,4).Value = 0

$nudos.Cells.Item(6,1).Value = 6
$nudos.Cells.Item(6,2).Value = 0
$nudos.Cells.Item(6,3).Value = 0
$nudos.Cells.Item(6,4).Value = 0

$nudos.Cells.Item(7,1).Value = 7
$nudos.Cells.Item(7,2).Value = 0
$nudos.Cells.Item(7,3).Value = 0
$nudos.Cells.Item(7,4).Value = 5000

$nudos.Cells.Item(8,1).Value = 8
$nudos.Cells.Item(8,2).Value = 6000
$nudos.Cells.Item(8,3).Value = 0
$nudos.Cells.Item(8,4).Value = 5000

# ---------------------------------------------------------------------
# Sheet "vxz" (local axis reference vector) - add nodes 4(value)-12
# ---------------------------------------------------------------------
$vxz = $wb.Worksheets.Item("vxz")

$vxz.Cells.Item(4,1).Value = 4
$vxz.Cells.Item(4,2).Value = 1
$vxz.Cells.Item(4,3).Value = 0
$vxz.Cells.Item(4,4).Value = 0

$vxz.Cells.Item(5,1).Value = 5
$vxz.Cells.Item(5,2).Value = 0
$vxz.Cells.Item(5,3).Value = 0
$vxz.Cells.Item(5,4).Value = 0

$vxz.Cells.Item(6,1).Value = 6
$vxz.Cells.Item(6,2).Value = 0
$vxz.Cells.Item(6,3).Value = 0
$vxz.Cells.Item(6,4).Value = 1

$vxz.Cells.Item(7,1).Value = 7
$vxz.Cells.Item(7,2).Value = 0
$vxz.Cells.Item(7,3).Value = 0
$vxz.Cells.Item(7,4).Value = 1

$vxz.Cells.Item(8,1).Value = 8
$vxz.Cells.Item(8,2).Value = 0
$vxz.Cells.Item(8,3).Value = 0
$vxz.Cells.Item(8,4).Value = 1

$vxz.Cells.Item(9,1).Value = 9
$vxz.Cells.Item(9,2).Value = 1
$vxz.Cells.Item(9,3).Value = 0
$vxz.Cells.Item(9,4).Value = 0

$vxz.Cells.Item(10,1).Value = 10
$vxz.Cells.Item(10,2).Value = 1
$vxz.Cells.Item(10,3).Value = 0
$vxz.Cells.Item(10,4).Value = 0

$vxz.Cells.Item(11,1).Value = 11
$vxz.Cells.Item(11,2).Value = 0
$vxz.Cells.Item(11,3).Value = 0
$vxz.Cells.Item(11,4).Value = 0

$vxz.Cells.Item(12,1).Value = 12
$vxz.Cells.Item(12,2).Value = 0
$vxz.Cells.Item(12,3).Value = 0
$vxz.Cells.Item(12,4).Value = 0

# ---------------------------------------------------------------------
# Sheet "conectividad" (element connectivity) - fix C1, add elements 4-12
# ---------------------------------------------------------------------
$conn = $wb.Worksheets.Item("conectividad")

$conn.Cells.Item(1,3).Value = 3

$conn.Cells.Item(4,1).Value = 4
$conn.Cells.Item(4,2).Value = 5
$conn.Cells.Item(4,3).Value = 1

$conn.Cells.Item(5,1).Value = 5
$conn.Cells.Item(5,2).Value = 6
$conn.Cells.Item(5,3).Value = 1

$conn.Cells.Item(6,1).Value = 6
$conn.Cells.Item(6,2).Value = 7
$conn.Cells.Item(6,3).Value = 1

$conn.Cells.Item(7,1).Value = 7
$conn.Cells.Item(7,2).Value = 7
$conn.Cells.Item(7,3).Value = 8

$conn.Cells.Item(8,1).Value = 8
$conn.Cells.Item(8,2).Value = 8
$conn.Cells.Item(8,3).Value = 3

$conn.Cells.Item(9,1).Value = 9
$conn.Cells.Item(9,2).Value = 6
$conn.Cells.Item(9,3).Value = 7

$conn.Cells.Item(10,1).Value = 10
$conn.Cells.Item(10,2).Value = 4
$conn.Cells.Item(10,3).Value = 8

$conn.Cells.Item(11,1).Value = 11
$conn.Cells.Item(11,2).Value = 4
$conn.Cells.Item(11,3).Value = 7

$conn.Cells.Item(12,1).Value = 12
$conn.Cells.Item(12,2).Value = 2
$conn.Cells.Item(12,3).Value = 1

# ---------------------------------------------------------------------
# Sheet "prop geom" (section properties) - replicate row 3 into rows 4-12
# so the "opensees" damage formulas (which read 'prop geom'!B9 / B11, etc.)
# stop evaluating to zero.
# ---------------------------------------------------------------------
$propGeom = $wb.Worksheets.Item("prop geom")
$mCol = [double]"7.8090000000000006E-9"

for ($r = 4; $r -le 12; $r++) {
    $propGeom.Cells.Item($r,1).Value  = $r
    $propGeom.Cells.Item($r,2).Value  = 11721.3
    $propGeom.Cells.Item($r,3).Value  = 120931398.7
    $propGeom.Cells.Item($r,4).Value  = 120931398.7
    $propGeom.Cells.Item($r,5).Value  = 241862797.40000001
    $propGeom.Cells.Item($r,6).Value  = 199947.98
    $propGeom.Cells.Item($r,7).Value  = 76903.070000000007
    $propGeom.Cells.Item($r,8).Value  = "circular"
    $propGeom.Cells.Item($r,9).Value  = "wo"
    $propGeom.Cells.Item($r,10).Value = 300
    $propGeom.Cells.Item($r,11).Value = 300
    $propGeom.Cells.Item($r,12).Value = 13
    $propGeom.Cells.Item($r,13).Value = $mCol
    $propGeom.Cells.Item($r,14).Value = 0.65
}

# ---------------------------------------------------------------------
# Sheet "fix nodes" - renumber node 1 as node 5, add supports for nodes 2,4,6
# (selected last so it keeps being the active/tabSelected sheet, matching
# the source file where it is the tab shown when the workbook is opened)
# ---------------------------------------------------------------------
$fixNodes = $wb.Worksheets.Item("fix nodes")

$fixNodes.Cells.Item(1,1).Value = 5

$fixNodes.Cells.Item(2,1).Value = 2
$fixNodes.Cells.Item(2,2).Value = 1
$fixNodes.Cells.Item(2,3).Value = 1
$fixNodes.Cells.Item(2,4).Value = 1
$fixNodes.Cells.Item(2,5).Value = 1
$fixNodes.Cells.Item(2,6).Value = 1
$fixNodes.Cells.Item(2,7).Value = 1

$fixNodes.Cells.Item(3,1).Value = 4
$fixNodes.Cells.Item(3,2).Value = 1
$fixNodes.Cells.Item(3,3).Value = 1
$fixNodes.Cells.Item(3,4).Value = 1
$fixNodes.Cells.Item(3,5).Value = 1
$fixNodes.Cells.Item(3,6).Value = 1
$fixNodes.Cells.Item(3,7).Value = 1

$fixNodes.Cells.Item(4,1).Value = 6
$fixNodes.Cells.Item(4,2).Value = 1
$fixNodes.Cells.Item(4,3).Value = 1
$fixNodes.Cells.Item(4,4).Value = 1
$fixNodes.Cells.Item(4,5).Value = 1
$fixNodes.Cells.Item(4,6).Value = 1
$fixNodes.Cells.Item(4,7).Value = 1

# ---------------------------------------------------------------------
# Update the on-screen selections to match the new, larger data ranges.
# "fix nodes" is selected last so it remains the active sheet/tab.
# ---------------------------------------------------------------------
$nudos.Range("A1:D8").Select()
$vxz.Range("A1:D12").Select()
$conn.Range("A1:C12").Select()
$propGeom.Range("A1:N12").Select()
$fixNodes.Range("A1:G4").Select()
